# GNA_metrics.xlsx update:
# Insert two new columns (sum_SASA, max_SASA) right after the existing
# "SASA" column (column C), shifting flexibility/Q/theta/conformation/
# monosaccharides/motifs/class two columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at D:E (pushes old D..J to F..L)
$ws.Columns("D:E").Insert()

# New header cells for the inserted columns
$ws.Cells.Item(1, 4).Value = "sum_SASA"
$ws.Cells.Item(1, 5).Value = "max_SASA"

# Make sure the new header cells carry the same header style as the rest
# of row 1 (bold, bordered, centered) by copying the style from C1.
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate sum_SASA (D) and max_SASA (E) values for each data row
$ws.Cells.Item(2, 4).Value = 5.076088470882665
$ws.Cells.Item(2, 5).Value = 2.753384233219247

$ws.Cells.Item(3, 4).Value = 4.502203803170425
$ws.Cells.Item(3, 5).Value = 2.280227683002259

$ws.Cells.Item(4, 4).Value = 2.651420443376771
$ws.Cells.Item(4, 5).Value = 2.651420443376771

$ws.Cells.Item(5, 4).Value = 2.351384937873273
$ws.Cells.Item(5, 5).Value = 2.351384937873273

$ws.Cells.Item(6, 4).Value = 2.270084840139852
$ws.Cells.Item(6, 5).Value = 2.270084840139852

$ws.Cells.Item(7, 4).Value = 4.737019782996375
$ws.Cells.Item(7, 5).Value = 2.456060678658576

$ws.Cells.Item(8, 4).Value = 2.423617038821429
$ws.Cells.Item(8, 5).Value = 2.423617038821429

$ws.Cells.Item(9, 4).Value = 7.478598491147868
$ws.Cells.Item(9, 5).Value = 2.575386834662778

$ws.Cells.Item(10, 4).Value = 5.32610747010928
$ws.Cells.Item(10, 5).Value = 2.76768229887049

$ws.Cells.Item(11, 4).Value = 5.053632587327314
$ws.Cells.Item(11, 5).Value = 2.530553721246874

$ws.Cells.Item(12, 4).Value = 5.224285506329593
$ws.Cells.Item(12, 5).Value = 2.637590683923242

$ws.Cells.Item(13, 4).Value = 4.610734076466935
$ws.Cells.Item(13, 5).Value = 2.415892868585273

$ws.Cells.Item(14, 4).Value = 5.020146705401829
$ws.Cells.Item(14, 5).Value = 2.545281005075958

$ws.Cells.Item(15, 4).Value = 4.58744599634927
$ws.Cells.Item(15, 5).Value = 2.381465013269206

$ws.Cells.Item(16, 4).Value = 4.622398868120259
$ws.Cells.Item(16, 5).Value = 2.325370683812271
